$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "F3745-TGGAACACAC",
    "F3746-ATCACTTGGA",
    "F3747-AAGTCAGGTT",
    "F3748-ACTAGGTCGA",
    "F3749-TCCAGTGACT",
    "F3750-GTAGAGACGA",
    "F3751-GCAGTGAGAT",
    "F3752-CTTCGATCTG",
    "F3753-GCATCAGTTC",
    "F3754-ACACAAGCAT",
    "F3755-ACACCTCTAC",
    "F3756-ACAAGTGTGT",
    "F3757-GGAACAGACT",
    "F3758-TAGTCTCGTC",
    "F3759-CAACCACCAT",
    "F3760-ATGAGACTTG",
    "F3761-ATCGACTCTT",
    "F3762-CGTCTACCTA",
    "F3763-AAGGTCAACC",
    "F3764-GCTGTACCTA",
    "F3765-TTGTCTGTAC",
    "F3766-CCTCTTGTTG",
    "F3767-AGCTAGAGTG",
    "F3768-AGAGTGAGTA",
    "F3769-AGCACACAGT",
    "F3770-CAGATCTGAT",
    "F3771-TCTACTACGA",
    "F3772-AACTGAGGTA",
    "F3773-TAGTCATCCT",
    "F3774-AGAGAAGTGC",
    "F3775-GAGACAACCT",
    "F3776-CGACTACTTG",
    "F3777-GATCTGACTC",
    "F3778-GGTGCAACAT",
    "F3779-ACAACCACGT",
    "F3780-TTCGATGCTG",
    "F3781-TTGAAGTCCA",
    "F3782-TGAGAGCTTC",
    "F3783-AGAACTCTGA",
    "F3784-GAACCTTGGT",
    "F3785-TCACAGAAGG",
    "F3786-TACGTAGTGA",
    "F3787-GGAAGAAGAT",
    "F3788-TCTAGCATGG",
    "F3789-CCATCTACAG",
    "F3790-TCTGTGTGTA",
    "F3791-AGACGATGTT",
    "F3792-GTCTACACCT",
    "F3793-TGACAGTCGA",
    "F3794-GAGTAGAGTC",
    "F3795-TACAAGTGCA",
    "F3796-AGGTCAGGAA",
    "F3797-AGTCACAAGA",
    "F3798-TCCACAAGGA",
    "F3799-TCAAGTCAAC",
    "F3800-AGAAGCAACC",
    "F3801-ATGTGAGTCG",
    "F3802-ACTTGACGAT",
    "F3803-TGTGATGGTC",
    "F3804-AGCAACCTTG",
    "F3805-ACGATCGTCT",
    "F3806-TCCTGAGATC",
    "F3807-CATCTGAGAC",
    "F3808-ATGAGGAGAG",
    "F3809-AGAGAGTACA",
    "F3810-TCCTAGATCC",
    "F3811-CAACATCTCC",
    "F3812-AGAAGGAAGT",
    "F3813-ACGTTCCATC",
    "F3814-CTTGAGTCAG",
    "F3815-ACTCTCTACA",
    "F3816-TGTGAAGAGT",
    "F3817-CGAAGGAGTT",
    "F3818-TCCAAGACTC",
    "F3819-AGTCTCGTTG",
    "F3820-ACTTCAGGAG",
    "F3821-GTCGAAGAGA",
    "F3822-CCTCAAGGAA",
    "F3823-CACTAGCTAC",
    "F3824-AGGTTCTGCA",
    "F3825-TCAGTAGAGA",
    "F3826-TCAGCAAGGT",
    "F3827-CGATCATCAT",
    "F3828-ATCCATCACC",
    "F3829-TAGGTTCGTA",
    "F3830-ATCTTGCACG",
    "F3831-CTCCTACTAG",
    "F3832-AAGTGCACGT",
    "F3833-ACTCACTGTG",
    "F3834-ATCCATGTAG",
    "F3835-GTTGTGAAGG",
    "F3836-AACAGAGAAC",
    "F3837-TCCACTCCTA",
    "F3838-TCGATGAAGT",
    "F3839-CGAAGAGAGT",
    "F3840-TACGTTCTCC"
)

$sequences = @(
    "AATGATACGGCGACCACCGAGATCTACACTGGAACACACTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACATCACTTGGATCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACAAGTCAGGTTTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACACTAGGTCGATCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTCCAGTGACTTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACGTAGAGACGATCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACGCAGTGAGATTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACCTTCGATCTGTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACGCATCAGTTCTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACACACAAGCATTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACACACCTCTACTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACACAAGTGTGTTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACGGAACAGACTTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTAGTCTCGTCTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACCAACCACCATTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACATGAGACTTGTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACATCGACTCTTTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACCGTCTACCTATCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACAAGGTCAACCTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACGCTGTACCTATCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTTGTCTGTACTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACCCTCTTGTTGTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACAGCTAGAGTGTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACAGAGTGAGTATCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACAGCACACAGTTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACCAGATCTGATTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTCTACTACGATCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACAACTGAGGTATCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTAGTCATCCTTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACAGAGAAGTGCTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACGAGACAACCTTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACCGACTACTTGTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACGATCTGACTCTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACGGTGCAACATTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACACAACCACGTTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTTCGATGCTGTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTTGAAGTCCATCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTGAGAGCTTCTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACAGAACTCTGATCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACGAACCTTGGTTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTCACAGAAGGTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTACGTAGTGATCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACGGAAGAAGATTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTCTAGCATGGTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACCCATCTACAGTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTCTGTGTGTATCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACAGACGATGTTTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACGTCTACACCTTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTGACAGTCGATCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACGAGTAGAGTCTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTACAAGTGCATCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACAGGTCAGGAATCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACAGTCACAAGATCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTCCACAAGGATCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTCAAGTCAACTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACAGAAGCAACCTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACATGTGAGTCGTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACACTTGACGATTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTGTGATGGTCTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACAGCAACCTTGTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACACGATCGTCTTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTCCTGAGATCTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACCATCTGAGACTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACATGAGGAGAGTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACAGAGAGTACATCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTCCTAGATCCTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACCAACATCTCCTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACAGAAGGAAGTTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACACGTTCCATCTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACCTTGAGTCAGTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACACTCTCTACATCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTGTGAAGAGTTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACCGAAGGAGTTTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTCCAAGACTCTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACAGTCTCGTTGTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACACTTCAGGAGTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACGTCGAAGAGATCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACCCTCAAGGAATCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACCACTAGCTACTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACAGGTTCTGCATCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTCAGTAGAGATCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTCAGCAAGGTTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACCGATCATCATTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACATCCATCACCTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTAGGTTCGTATCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACATCTTGCACGTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACCTCCTACTAGTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACAAGTGCACGTTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACACTCACTGTGTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACATCCATGTAGTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACGTTGTGAAGGTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACAACAGAGAACTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTCCACTCCTATCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTCGATGAAGTTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACCGAAGAGAGTTCGTCGGCAGCGTC",
    "AATGATACGGCGACCACCGAGATCTACACTACGTTCTCCTCGTCGGCAGCGTC"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $names[$i]
    $ws.Cells.Item($row, 3).Value = $sequences[$i]
}
